# Add "Giornata 8" data column and recompute "Total average" across 8 giornate
# for FantaMedia_Squadre sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First, duplicate the header style from I1 into J1 (copies value+format),
# then overwrite the values so the header order becomes:
# ... H1=Giornata 7, I1=Giornata 8, J1=Total average
$ws.Range("I1").Copy($ws.Range("J1"))
$ws.Range("J1").Value = "Total average "
$ws.Range("I1").Value = "Giornata 8"

# New Giornata 8 score per team (column I), and the recomputed
# Total average across Giornata 1-8 (column J).
$ws.Range("I2").Value = 5.785714285714286
$ws.Range("J2").Value = 5.976219093406593

$ws.Range("I3").Value = 6.615384615384615
$ws.Range("J3").Value = 5.872856570512821

$ws.Range("I4").Value = 6.576923076923077
$ws.Range("J4").Value = 6.033012820512822

$ws.Range("I5").Value = 5.5
$ws.Range("J5").Value = 5.984775641025641

$ws.Range("I6").Value = 5.692307692307693
$ws.Range("J6").Value = 5.804601648351649

$ws.Range("I7").Value = 5.333333333333333
$ws.Range("J7").Value = 5.834612262737263

$ws.Range("I8").Value = 5.642857142857143
$ws.Range("J8").Value = 5.772744963369964

$ws.Range("I9").Value = 6.066666666666666
$ws.Range("J9").Value = 5.944432773109243

$ws.Range("I10").Value = 6.576923076923077
$ws.Range("J10").Value = 6.175730519480521

$ws.Range("I11").Value = 6.428571428571429
$ws.Range("J11").Value = 5.995386904761904

$ws.Range("I12").Value = 6.730769230769231
$ws.Range("J12").Value = 6.317055860805861

$ws.Range("I13").Value = 5.730769230769231
$ws.Range("J13").Value = 6.206267690642691

$ws.Range("I14").Value = 5.071428571428571
$ws.Range("J14").Value = 5.775869963369964

$ws.Range("I15").Value = 6.833333333333333
$ws.Range("J15").Value = 6.327953296703297

$ws.Range("I16").Value = 5.576923076923077
$ws.Range("J16").Value = 6.009354967948717

$ws.Range("I17").Value = 6.0625
$ws.Range("J17").Value = 6.110857371794872

$ws.Range("I18").Value = 6.166666666666667
$ws.Range("J18").Value = 5.908482142857142

$ws.Range("I19").Value = 6
$ws.Range("J19").Value = 5.793326465201465

$ws.Range("I20").Value = 6.307692307692307
$ws.Range("J20").Value = 5.962912087912088

$ws.Range("I21").Value = 6.033333333333333
$ws.Range("J21").Value = 6.120089285714285
